$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 68.158272
$ws.Range("H2").Value = 204.474816
$ws.Range("I2").Value = 0.164824640128582
$ws.Range("J2").Value = 0.1648246401285819
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.051093
$ws.Range("N2").Value = 0.153279
$ws.Range("O2").Value = 0.01450579975525089
$ws.Range("P2").Value = 0.01450579975525089
$ws.Range("Q2").Value = 3.482410591296
$ws.Range("R2").Value = 31.341695321664
$ws.Range("S2").Value = 0.0023909132244365
$ws.Range("T2").Value = 0.0023909132244365

$ws.Range("G3").Value = 68.158272
$ws.Range("H3").Value = 204.474816
$ws.Range("I3").Value = 0.164824640128582
$ws.Range("J3").Value = 0.1648246401285819
$ws.Range("O3").Value = 0.2313022967634575
$ws.Range("P3").Value = 0.2313022967634575
$ws.Range("Q3").Value = 55.52879411206399
$ws.Range("R3").Value = 499.759147008576
$ws.Range("S3").Value = 0.03812431782495135
$ws.Range("T3").Value = 0.03812431782495135

$ws.Range("G4").Value = 68.158272
$ws.Range("H4").Value = 204.474816
$ws.Range("I4").Value = 0.164824640128582
$ws.Range("J4").Value = 0.1648246401285819
$ws.Range("M4").Value = 2.656449666666667
$ws.Range("N4").Value = 7.969348999999999
$ws.Range("O4").Value = 0.7541919034812916
$ws.Range("P4").Value = 0.7541919034812917
$ws.Range("Q4").Value = 181.059018934976
$ws.Range("R4").Value = 1629.531170414784
$ws.Range("S4").Value = 0.1243094090791941
$ws.Range("T4").Value = 0.1243094090791941

$ws.Range("I5").Value = 0.3471155005059974
$ws.Range("J5").Value = 0.3471155005059974
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.051093
$ws.Range("N5").Value = 0.153279
$ws.Range("O5").Value = 0.01450579975525089
$ws.Range("P5").Value = 0.01450579975525089
$ws.Range("Q5").Value = 7.333847017182
$ws.Range("R5").Value = 66.004623154638
$ws.Range("S5").Value = 0.005035187942283686
$ws.Range("T5").Value = 0.005035187942283687

$ws.Range("I6").Value = 0.3471155005059974
$ws.Range("J6").Value = 0.3471155005059974
$ws.Range("O6").Value = 0.2313022967634575
$ws.Range("P6").Value = 0.2313022967634575
$ws.Range("S6").Value = 0.08028861250923429
$ws.Range("T6").Value = 0.08028861250923429

$ws.Range("I7").Value = 0.3471155005059974
$ws.Range("J7").Value = 0.3471155005059974
$ws.Range("M7").Value = 2.656449666666667
$ws.Range("N7").Value = 7.969348999999999
$ws.Range("O7").Value = 0.7541919034812916
$ws.Range("P7").Value = 0.7541919034812917
$ws.Range("Q7").Value = 381.3045909259087
$ws.Range("R7").Value = 3431.741318333178
$ws.Range("S7").Value = 0.2617917000544794
$ws.Range("T7").Value = 0.2617917000544794

$ws.Range("G8").Value = 201.822474
$ws.Range("H8").Value = 605.4674219999999
$ws.Range("I8").Value = 0.4880598593654206
$ws.Range("J8").Value = 0.4880598593654206
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.051093
$ws.Range("N8").Value = 0.153279
$ws.Range("O8").Value = 0.01450579975525089
$ws.Range("P8").Value = 0.01450579975525089
$ws.Range("Q8").Value = 10.311715664082
$ws.Range("R8").Value = 92.80544097673798
$ws.Range("S8").Value = 0.0070796985885307
$ws.Range("T8").Value = 0.007079698588530701

$ws.Range("G9").Value = 201.822474
$ws.Range("H9").Value = 605.4674219999999
$ws.Range("I9").Value = 0.4880598593654206
$ws.Range("J9").Value = 0.4880598593654206
$ws.Range("O9").Value = 0.2313022967634575
$ws.Range("P9").Value = 0.2313022967634575
$ws.Range("Q9").Value = 164.425509583538
$ws.Range("R9").Value = 1479.829586251842
$ws.Range("S9").Value = 0.1128893664292718
$ws.Range("T9").Value = 0.1128893664292718

$ws.Range("G10").Value = 201.822474
$ws.Range("H10").Value = 605.4674219999999
$ws.Range("I10").Value = 0.4880598593654206
$ws.Range("J10").Value = 0.4880598593654206
$ws.Range("M10").Value = 2.656449666666667
$ws.Range("N10").Value = 7.969348999999999
$ws.Range("O10").Value = 0.7541919034812916
$ws.Range("P10").Value = 0.7541919034812917
$ws.Range("Q10").Value = 536.1312437831419
$ws.Range("R10").Value = 4825.181194048278
$ws.Range("S10").Value = 0.368090794347618
$ws.Range("T10").Value = 0.3680907943476181
